$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stage the full target B2:E51 block as text in a scratch area (columns Z:AC),
# then copy+paste-special (values only) onto the live range so the destination
# keeps its original (default) cell style while its value becomes text - matching
# how Excel stored these numeric-looking strings as inline text in the source file.
$scratch = $ws.Range("Z2:AC51")
$scratch.NumberFormat = "@"

$ws.Range("Z2").Value = 'Bitcoin'
$ws.Range("AA2").Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range("AB2").Value = '68.063.59'
$ws.Range("AC2").Value = '  +1.38%  '
$ws.Range("Z3").Value = 'Ethereum'
$ws.Range("AA3").Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range("AB3").Value = '3.530.72'
$ws.Range("AC3").Value = '  +0.34%  '
$ws.Range("Z4").Value = 'TetherUSD'
$ws.Range("AA4").Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range("AB4").Value = '1.00'
$ws.Range("AC4").Value = '  -0.01%  '
$ws.Range("Z5").Value = 'BNB'
$ws.Range("AA5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("AB5").Value = '601.15'
$ws.Range("AC5").Value = '  +1.14%  '
$ws.Range("Z6").Value = 'Solana'
$ws.Range("AA6").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("AB6").Value = '183.69'
$ws.Range("AC6").Value = '  +5.56%  '
$ws.Range("Z7").Value = 'USDC'
$ws.Range("AA7").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("AB7").Value = '1.00'
$ws.Range("AC7").Value = '  +0.02%  '
$ws.Range("Z8").Value = 'XRP'
$ws.Range("AA8").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("AB8").Value = '0.599'
$ws.Range("AC8").Value = '  +0.52%  '
$ws.Range("Z9").Value = 'Dogecoin'
$ws.Range("AA9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("AB9").Value = '0.140'
$ws.Range("AC9").Value = '  +4.64%  '
$ws.Range("Z10").Value = 'Toncoin'
$ws.Range("AA10").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("AB10").Value = '7.15'
$ws.Range("AC10").Value = '  -1.77%  '
$ws.Range("Z11").Value = 'Cardano'
$ws.Range("AA11").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("AB11").Value = '0.443'
$ws.Range("AC11").Value = '  +1.30%  '
$ws.Range("Z12").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("AA12").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("AB12").Value = '4.143.48'
$ws.Range("AC12").Value = '  +0.36%  '
$ws.Range("Z13").Value = 'Avalanche'
$ws.Range("AA13").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("AB13").Value = '32.58'
$ws.Range("AC13").Value = '  +12.18%  '
$ws.Range("Z14").Value = 'TRON'
$ws.Range("AA14").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("AB14").Value = '0.134'
$ws.Range("AC14").Value = '  -0.15%  '
$ws.Range("Z15").Value = 'WrappedBTC'
$ws.Range("AA15").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("AB15").Value = '68.018.29'
$ws.Range("AC15").Value = '  +1.36%  '
$ws.Range("Z16").Value = 'ShibaInu'
$ws.Range("AA16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("AB16").Value = '0.0000183'
$ws.Range("AC16").Value = '  +0.82%  '
$ws.Range("Z17").Value = 'WrappedEther'
$ws.Range("AA17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("AB17").Value = '3.535.20'
$ws.Range("AC17").Value = '  +0.35%  '
$ws.Range("Z18").Value = 'Polkadot'
$ws.Range("AA18").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("AB18").Value = '6.42'
$ws.Range("AC18").Value = '  +1.39%  '
$ws.Range("Z19").Value = 'Chainlink'
$ws.Range("AA19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("AB19").Value = '15.00'
$ws.Range("AC19").Value = '  +4.92%  '
$ws.Range("Z20").Value = 'BitcoinCash'
$ws.Range("AA20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("AB20").Value = '399.66'
$ws.Range("AC20").Value = '  +0.61%  '
$ws.Range("Z21").Value = 'Uniswap'
$ws.Range("AA21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("AB21").Value = '8.14'
$ws.Range("AC21").Value = '  +1.79%  '
$ws.Range("Z22").Value = 'Litecoin'
$ws.Range("AA22").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("AB22").Value = '73.62'
$ws.Range("AC22").Value = '  +0.35%  '
$ws.Range("Z23").Value = 'Polygon'
$ws.Range("AA23").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("AB23").Value = '0.548'
$ws.Range("AC23").Value = '  +1.14%  '
$ws.Range("Z24").Value = 'Dai'
$ws.Range("AA24").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("AB24").Value = '1.00'
$ws.Range("AC24").Value = '  +0.06%  '
$ws.Range("Z25").Value = 'PEPE'
$ws.Range("AA25").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("AB25").Value = '0.0000126'
$ws.Range("AC25").Value = '  +2.81%  '
$ws.Range("Z26").Value = 'LEO'
$ws.Range("AA26").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("AB26").Value = '5.71'
$ws.Range("AC26").Value = '  +0.20%  '
$ws.Range("Z27").Value = 'InternetComputer(DFINITY)'
$ws.Range("AA27").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("AB27").Value = '10.79'
$ws.Range("AC27").Value = '  +5.59%  '
$ws.Range("Z28").Value = 'Kaspa'
$ws.Range("AA28").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("AB28").Value = '0.179'
$ws.Range("AC28").Value = '  -0.98%  '
$ws.Range("Z29").Value = 'Binance-PegBSC-USD'
$ws.Range("AA29").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("AB29").Value = '0.997'
$ws.Range("AC29").Value = '  -0.10%  '
$ws.Range("Z30").Value = 'NEARProtocol'
$ws.Range("AA30").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("AB30").Value = '6.32'
$ws.Range("AC30").Value = '  +0.52%  '
$ws.Range("Z31").Value = 'Fetch.AI'
$ws.Range("AA31").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("AB31").Value = '1.47'
$ws.Range("AC31").Value = '  +1.24%  '
$ws.Range("Z32").Value = 'PancakeSwap'
$ws.Range("AA32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("AB32").Value = '2.09'
$ws.Range("AC32").Value = '  +0.82%  '
$ws.Range("Z33").Value = 'EthereumClassic'
$ws.Range("AA33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("AB33").Value = '24.15'
$ws.Range("AC33").Value = '  +0.51%  '
$ws.Range("Z34").Value = 'Aptos'
$ws.Range("AA34").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("AB34").Value = '7.52'
$ws.Range("AC34").Value = '  +1.36%  '
$ws.Range("Z35").Value = 'USDe'
$ws.Range("AA35").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("AB35").Value = '1.00'
$ws.Range("AC35").Value = '  +0.13%  '
$ws.Range("Z36").Value = 'ImmutableX'
$ws.Range("AA36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("AB36").Value = '1.69'
$ws.Range("AC36").Value = '  +2.79%  '
$ws.Range("Z37").Value = 'Monero'
$ws.Range("AA37").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("AB37").Value = '164.04'
$ws.Range("AC37").Value = '  +0.21%  '
$ws.Range("Z38").Value = 'Stacks'
$ws.Range("AA38").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("AB38").Value = '1.98'
$ws.Range("AC38").Value = '  +3.39%  '
$ws.Range("Z39").Value = 'Mantle'
$ws.Range("AA39").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("AB39").Value = '0.881'
$ws.Range("AC39").Value = '  -1.49%  '
$ws.Range("Z40").Value = 'RenderToken'
$ws.Range("AA40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("AB40").Value = '7.19'
$ws.Range("AC40").Value = '  +4.01%  '
$ws.Range("Z41").Value = 'dogwifhat'
$ws.Range("AA41").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("AB41").Value = '2.79'
$ws.Range("AC41").Value = '  +6.95%  '
$ws.Range("Z42").Value = 'Filecoin'
$ws.Range("AA42").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("AB42").Value = '4.78'
$ws.Range("AC42").Value = '  +1.64%  '
$ws.Range("Z43").Value = 'EnergySwap'
$ws.Range("AA43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("AB43").Value = '27.16'
$ws.Range("AC43").Value = '  +2.52%  '
$ws.Range("Z44").Value = 'InjectiveProtocol'
$ws.Range("AA44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("AB44").Value = '27.74'
$ws.Range("AC44").Value = '  -0.32%  '
$ws.Range("Z45").Value = 'Maker'
$ws.Range("AA45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("AB45").Value = '2.888.58'
$ws.Range("AC45").Value = '  +3.03%  '
$ws.Range("Z46").Value = 'Hedera'
$ws.Range("AA46").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("AB46").Value = '0.0743'
$ws.Range("AC46").Value = '  -0.46%  '
$ws.Range("Z47").Value = 'OKB'
$ws.Range("AA47").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("AB47").Value = '42.42'
$ws.Range("AC47").Value = '  -1.02%  '
$ws.Range("Z48").Value = 'Bittensor'
$ws.Range("AA48").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("AB48").Value = '353.49'
$ws.Range("AC48").Value = '  +4.04%  '
$ws.Range("Z49").Value = 'VeChain'
$ws.Range("AA49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("AB49").Value = '0.0307'
$ws.Range("AC49").Value = '  +0.59%  '
$ws.Range("Z50").Value = 'ONDO'
$ws.Range("AA50").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("AB50").Value = '1.09'
$ws.Range("AC50").Value = '  -0.70%  '
$ws.Range("Z51").Value = 'Arweave'
$ws.Range("AA51").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("AB51").Value = '34.17'
$ws.Range("AC51").Value = '  +1.99%  '

$scratch.Copy()
$ws.Range("B2:E51").PasteSpecial(-4163)
$ws.Application.CutCopyMode = $false
$scratch.Clear()
